$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14. This pushes the existing rows 14-51
# down to 15-52, preserving all of their data (including the row-14
# formatting carried onto the new blank row's D column).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record's data.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44497
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112021
$ws.Range("G14").Value = "Ají"
$ws.Range("H14").Value = "Americana (o)"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 14
$ws.Range("K14").Value = 36000
$ws.Range("L14").Value = 37000
$ws.Range("M14").Value = 36571
$ws.Range("N14").Value = "$/caja 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 1463
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
